$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old shuffled-address rows (rows 2-6), keep row 1 formatting intact
$ws.Rows("2:6").Delete()

# Set the single remaining value - the canonical parsed address string
$ws.Range("A1").Value = "107 OLD COLONY RD    RICHMOND HILL ON,         L4E3X2    "

# Reset selection to A1 (previous selection pointed at J5)
$ws.Range("A1").Select()
